$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.703.22"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.644.90"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.92"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.38"
$ws.Range("E8").Value = "  +1.37%  "
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.877.82"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.649.30"
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.77"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.690.10"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.36"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.70"
$ws.Range("E19").Value = "  +1.87%  "
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("E23").Value = "  +6.07%  "
$ws.Range("E24").Value = "  -2.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.10"
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  -1.05%  "
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0488"
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("E32").Value = "  +1.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.454.02"
$ws.Range("E33").Value = "  +3.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.14"
$ws.Range("E34").Value = "  +0.69%  "
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.570"
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.885"
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("E39").Value = "  +0.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.887"
$ws.Range("E40").Value = "  +12.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "71.15"
$ws.Range("E41").Value = "  +10.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.03"
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("E44").Value = "  +2.52%  "
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("E46").Value = "  +0.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.787.56"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.73"
$ws.Range("E48").Value = "  +5.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.84"
$ws.Range("E49").Value = "  -1.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0991"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.77"
$ws.Range("E51").Value = "  +0.88%  "
